$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 160, shifting existing rows
# 160-165 down to 162-167 (the rest of the weekly "Ají" price records).
$ws.Rows.Item(160).Insert()
$ws.Rows.Item(160).Insert()

# New row 160: Americana (o) / Primera, Provincia de Limarí
$ws.Cells.Item(160, 1).Value = 8
$ws.Cells.Item(160, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44568
$ws.Cells.Item(160, 5).Value = 4
$ws.Cells.Item(160, 6).Value = 100112021
$ws.Cells.Item(160, 7).Value = "Ají"
$ws.Cells.Item(160, 8).Value = "Americana (o)"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 600
$ws.Cells.Item(160, 11).Value = 24000
$ws.Cells.Item(160, 12).Value = 25000
$ws.Cells.Item(160, 13).Value = 24500
$ws.Cells.Item(160, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(160, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(160, 16).Value = 980
$ws.Cells.Item(160, 17).Value = 25
$ws.Cells.Item(160, 18).Value = "Hortaliza"
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(162, 4).NumberFormat

# New row 161: Inferno / Primera, Provincia de Limarí
$ws.Cells.Item(161, 1).Value = 8
$ws.Cells.Item(161, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(161, 3).Value = "Coquimbo"
$ws.Cells.Item(161, 4).Value = 44568
$ws.Cells.Item(161, 5).Value = 4
$ws.Cells.Item(161, 6).Value = 100112021
$ws.Cells.Item(161, 7).Value = "Ají"
$ws.Cells.Item(161, 8).Value = "Inferno"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 600
$ws.Cells.Item(161, 11).Value = 14000
$ws.Cells.Item(161, 12).Value = 15000
$ws.Cells.Item(161, 13).Value = 14500
$ws.Cells.Item(161, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(161, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(161, 16).Value = 967
$ws.Cells.Item(161, 17).Value = 15
$ws.Cells.Item(161, 18).Value = "Hortaliza"
$ws.Cells.Item(161, 4).NumberFormat = $ws.Cells.Item(162, 4).NumberFormat
